$d = $word.ActiveDocument

$pairs = @(
    @("2024-09-11 Wednesday", "2024-09-12 Thursday"),
    @("61-55=", "0+11="),
    @("80-5=", "72-27="),
    @("7+32=", "73-1="),
    @("84-12=", "74-14="),
    @("33+18=", "95-26="),
    @("63-0=", "21+57="),
    @("28+54=", "99-55="),
    @("22-19=", "49+24="),
    @("95-16=", "46+6="),
    @("68+1=", "30+19="),
    @("46-43=", "29+22="),
    @("1+94=", "48+32="),
    @("87-6=", "16+41="),
    @("80-32=", "49+5="),
    @("42+21=", "16+8="),
    @("72-29=", "96-8="),
    @("81-61=", "56-14="),
    @("62+4=", "55+26="),
    @("90-76=", "40-3="),
    @("21-15=", "67-15="),
    @("55+0=", "9+78="),
    @("80-11=", "96-36="),
    @("95+3=", "73-69="),
    @("1+65=", "32+54="),
    @("13+10=", "61-59="),
    @("1+51=", "36+41="),
    @("63-11=", "82-13="),
    @("33+10=", "61-36="),
    @("73+22=", "22+72="),
    @("33-15=", "26+44="),
    @("12-2=", "26+8="),
    @("49+23=", "67-4="),
    @("16-10=", "50-47="),
    @("30+49=", "92-9="),
    @("24+35=", "72-16="),
    @("99-21=", "38+35="),
    @("49+12=", "85-32="),
    @("4+1=", "37+13="),
    @("98-25=", "38+36="),
    @("77+17=", "39+32="),
    @("39+50=", "47+18="),
    @("71+3=", "57-13="),
    @("62-29=", "86-40="),
    @("90+1=", "11+49="),
    @("52+38=", "14+81="),
    @("62-42=", "61-31="),
    @("98-88=", "48+37="),
    @("13-7=", "0+31="),
    @("71-5=", "47+31="),
    @("37-29=", "54-32="),
    @("65+21=", "75-49="),
    @("52-10=", "42+19="),
    @("99-2=", "87-52="),
    @("46-23=", "80-41="),
    @("11+23=", "27-16="),
    @("41-20=", "87-72="),
    @("0+7=", "6-2="),
    @("32+10=", "60-37="),
    @("17+7=", "48+44="),
    @("12+53=", "36+46="),
    @("59-46=", "44+43="),
    @("42+37=", "8+12="),
    @("9+16=", "6+29="),
    @("30-23=", "26-10="),
    @("25+40=", "84-60="),
    @("99-58=", "65-17="),
    @("22+31=", "4+85="),
    @("90-64=", "84+6="),
    @("17+75=", "8+11="),
    @("43+0=", "72-65="),
    @("94-89=", "49+42="),
    @("64-61=", "8+41="),
    @("65-33=", "43+52="),
    @("18+4=", "8-1="),
    @("98-26=", "75-27="),
    @("46-9=", "34+40="),
    @("96-10=", "14-0="),
    @("66-33=", "96-40="),
    @("86-64=", "99-95="),
    @("92-6=", "48+42="),
    @("91-87=", "36+27="),
    @("91-62=", "6+27="),
    @("49+21=", "87-49="),
    @("77-45=", "87+9="),
    @("84-62=", "78-45="),
    @("80-36=", "31+42="),
    @("93-71=", "31+49="),
    @("59+17=", "12+21="),
    @("60-14=", "0+84="),
    @("17-1=", "48+46="),
    @("84+14=", "71-43="),
    @("97-55=", "80-73="),
    @("84-45=", "75-39="),
    @("51-4=", "77-43="),
    @("67-16=", "48+0="),
    @("37-37=", "46-2="),
    @("34-23=", "3+88="),
    @("73-20=", "90-75="),
    @("89+1=", "8+78="),
    @("44+55=", "35+19="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Replaced $($pairs.Count) items"
